$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Developer name
$ws.Range("C3").Value = "Damien Altenburg"

# Row 7 - __init__ / Attribute set to input value
$ws.Range("E7").Value = "None"
$ws.Range("F7").Value = "name = ""ISD""" + "`n" + "department = Department.COMPUTER_SCIENCE" + "`n" + "credit_hours = 90" + "`n" + "capacity = 30" + "`n" + "current_enrollment = 14" + "`n" + "lab_eqiupment = ""Laptop"""
$ws.Range("G7").Value = "Object initialized with the correct state"

# Row 8 - __init__ / Exception raised when lab equipment is blank
$ws.Range("E8").Value = "None"
$ws.Range("F8").Value = "name = ""ISD""" + "`n" + "department = Department.COMPUTER_SCIENCE" + "`n" + "credit_hours = 90" + "`n" + "capacity = 30" + "`n" + "current_enrollment = 14" + "`n" + "lab_eqiupment = """""
$ws.Range("G8").Value = "Object initialized with the correct state (lab_equipment = ""None"")"

# Row 9 - __str__ / returns string in expected format
$ws.Range("E9").Value = "name = ""ISD""" + "`n" + "department = Department.COMPUTER_SCIENCE" + "`n" + "credit_hours = 90" + "`n" + "capacity = 30" + "`n" + "current_enrollment = 14" + "`n" + "lab_eqiupment = ""Laptop"""
$ws.Range("F9").Value = "None"
$ws.Range("G9").Value = "Course: ISD" + "`n" + "Department: Computer Science" + "`n" + "Credit Hours: 90" + "`n" + "Lab Equipment: Laptop"

# Row 10 - enroll_student / Successfully enrolls a student when there is sufficient capacity.
$ws.Range("E10").Value = "name = ""ISD""" + "`n" + "department = Department.COMPUTER_SCIENCE" + "`n" + "credit_hours = 90" + "`n" + "capacity = 30" + "`n" + "current_enrollment = 14" + "`n" + "lab_eqiupment = ""Laptop"""
$ws.Range("F10").Value = "Student(student_number = 123, name  = ""Joe Smith"", department = Department.COMPUTER_SCIENCE)"
$ws.Range("G10").Value = """Joe Smith has been successfully enrolled in ISD."""

# Row 11 - enroll_student / Rejects enrollment when there is insufficient capacity.
$ws.Range("E11").Value = "name = ""ISD""" + "`n" + "department = Department.COMPUTER_SCIENCE" + "`n" + "credit_hours = 90" + "`n" + "capacity = 15" + "`n" + "current_enrollment = 15" + "`n" + "lab_eqiupment = ""Laptop"""
$ws.Range("F11").Value = "Student(student_number = 123, name  = ""Joe Smith"", department = Department.COMPUTER_SCIENCE)"
$ws.Range("G11").Value = """Joe Smith has NOT been enrolled in lab: ISD due to insufficient capacity."""

# These cells move from the plain wrap-text style to the bold wrap-text style
$ws.Range("F9").Font.Bold = $true
$ws.Range("E10:G10").Font.Bold = $true
$ws.Range("E11:G11").Font.Bold = $true

# Row heights - set explicit height (auto row height became 105 for all wrapped rows)
$ws.Rows.Item(7).RowHeight = 105
$ws.Rows.Item(8).RowHeight = 105
$ws.Rows.Item(9).RowHeight = 105
$ws.Rows.Item(10).RowHeight = 105
$ws.Rows.Item(11).RowHeight = 105

# Column widths
$ws.Columns.Item(5).ColumnWidth = 40.5703125
$ws.Columns.Item(6).ColumnWidth = 42.85546875
$ws.Columns.Item(7).ColumnWidth = 56.7109375

# Selection / view
$ws.Range("F11").Select()
$excel.ActiveWindow.ScrollRow = 5
